# Regenerate the handback-status report for the two e2e files that were
# re-processed: the old "acb3d08b..." file became "32fa0c8e..." and the
# old "ae20c40a..." file became "fffffc8c9a4e...".

$wb = $excel.ActiveWorkbook

$oldGuid1 = "acb3d08b-601e-4505-b3a3-5b94ba208151"
$newGuid1 = "32fa0c8e-2890-4fb7-90f9-7212dd3922ef"
$oldGuid2 = "ae20c40a-579e-4708-88c6-9d041cfce420"
$newGuid2 = "fffffc8c9a4e-3817-4304-8563-e63318d4b77d"

$newMd1 = "$newGuid1.md"
$newMd2 = "$newGuid2.md"
$newMd1Disp = "e2e\$newGuid1.md"
$newMd2Disp = "e2e\$newGuid2.md"

$newXlfZh = "$newGuid1.1b21172b5759f6b658c312369ecda99410f0a9d4.zh-cn.xlf"
$newXlfDe = "$newGuid1.1b21172b5759f6b658c312369ecda99410f0a9d4.de-de.xlf"

$newHoDate = "2016-09-05 05:09:10"
$newZhHoffDate = "2016-09-05 05:09:01"
$newZhHbackDate = "2016-09-05 05:09:29"
$newDeHbackDate = "2016-09-05 05:09:37"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newMd1
$wsOverview.Range("A3").Value = $newMd2
$wsOverview.Range("G2").Value = $newHoDate
$wsOverview.Range("G3").Value = $newHoDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/60b9f927591dcf76bb6b6345f2dc81006378ac7d/e2e/$oldGuid1.md", [Type]::Missing, [Type]::Missing, $newMd1Disp)
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/60b9f927591dcf76bb6b6345f2dc81006378ac7d/e2e/$oldGuid2.md", [Type]::Missing, [Type]::Missing, $newMd2Disp)

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("G2").Value = $newXlfZh
$wsZh.Range("H2").Value = $newZhHoffDate
$wsZh.Range("J2").Value = $newXlfZh
$wsZh.Range("K2").Value = $newZhHbackDate

$wsZh.Range("G3").Value = $newXlfZh
$wsZh.Range("H3").Value = $newZhHoffDate
$wsZh.Range("J3").Value = $newXlfZh
$wsZh.Range("K3").Value = $newZhHbackDate

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/60b9f927591dcf76bb6b6345f2dc81006378ac7d/e2e/$oldGuid1.md", [Type]::Missing, [Type]::Missing, $newMd1)
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/4165315289e9363b44cc2821beda224979a248e4/e2e/$oldGuid1.md", [Type]::Missing, [Type]::Missing, $newMd1)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/60b9f927591dcf76bb6b6345f2dc81006378ac7d/e2e/$oldGuid2.md", [Type]::Missing, [Type]::Missing, $newMd2)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/4165315289e9363b44cc2821beda224979a248e4/e2e/$oldGuid2.md", [Type]::Missing, [Type]::Missing, $newMd2)

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("G2").Value = $newXlfDe
$wsDe.Range("H2").Value = $newHoDate
$wsDe.Range("J2").Value = $newXlfDe
$wsDe.Range("K2").Value = $newDeHbackDate

$wsDe.Range("G3").Value = $newXlfDe
$wsDe.Range("H3").Value = $newHoDate
$wsDe.Range("J3").Value = $newXlfDe
$wsDe.Range("K3").Value = $newDeHbackDate

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/60b9f927591dcf76bb6b6345f2dc81006378ac7d/e2e/$oldGuid1.md", [Type]::Missing, [Type]::Missing, $newMd1)
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/8888e362fe93e401eed71ad4ebc899d31a112774/e2e/$oldGuid1.md", [Type]::Missing, [Type]::Missing, $newMd1)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/60b9f927591dcf76bb6b6345f2dc81006378ac7d/e2e/$oldGuid2.md", [Type]::Missing, [Type]::Missing, $newMd2)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/8888e362fe93e401eed71ad4ebc899d31a112774/e2e/$oldGuid2.md", [Type]::Missing, [Type]::Missing, $newMd2)
